# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (column F) counts on the "展览" sheet (rows 2-36)
# and on the "全部类型" sheet (rows 3-38, same data shifted down by one row).

$wb = $excel.ActiveWorkbook

# Map of row -> new F value for the "展览" worksheet
$exhibitionUpdates = @{
    2  = 3129
    3  = 528
    4  = 1089
    5  = 82
    9  = 1122
    10 = 15600
    12 = 170
    13 = 1023
    14 = 6151
    19 = 112
    21 = 28
    23 = 10
    27 = 862
    29 = 4993
    30 = 479
    31 = 11024
    35 = 159
    36 = 3795
}

# Map of row -> new F value for the "全部类型" worksheet (same records, offset by 1 row)
$allTypesUpdates = @{
    3  = 3129
    4  = 528
    5  = 1089
    6  = 82
    10 = 1122
    11 = 15600
    13 = 170
    14 = 1023
    15 = 6151
    20 = 112
    22 = 28
    24 = 10
    28 = 862
    30 = 4993
    31 = 479
    33 = 11024
    37 = 159
    38 = 3795
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
